$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new quote row (row 5). Cells are written quote -> author -> topic to
# match the shared-string insertion order of the source edit.
$ws.Range("C5").Value = "Really successful people feel the same lack of motivation as everyone else.  The difference is that they still find a way to show up despite the feelings of boredom."
$ws.Range("B5").Value = "James Clear"
$ws.Range("A5").Value = "habits"

# Add the new quote row (row 6). Quote -> topic -> author ordering, matching
# the source edit's shared-string insertion order.
$ws.Range("C6").Value = "Goals are good for setting a direction, but systems are best for making progress."
$ws.Range("A6").Value = "systems"
$ws.Range("B6").Value = "James Clear"

# Match the row heights Excel would have produced for the wrapped text in
# column C (two text lines for the long quote in row 5, one line in row 6) --
# the same heights already used by the sheet's other rows (ht="34"/"17").
$ws.Rows.Item(5).RowHeight = 34
$ws.Rows.Item(6).RowHeight = 17

# Update the selection to reflect the new first empty row, as in the diff
$ws.Range("A7").Select() | Out-Null
